# Populate the "testData" sheet with the login/user test data table,
# add a mailto hyperlink on the e-mail cell, and add a new, initially
# empty "writeData" worksheet right after it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 1 (+ QA contact e-mail in E1)
$ws.Range("A1").Value = "rahul"
$ws.Range("B1").Value = "ritesh"
$ws.Range("C1").Value = 123
$ws.Range("D1").Value = "tester"
$ws.Range("E1").Value = "rahulriteshqa@gmail.com"

# Row 2
$ws.Range("B2").Value = "admin123"
$ws.Range("A2").Value = "Admin"
$ws.Range("C2").Value = 124

# Row 3
$ws.Range("A3").Value = "username1"
$ws.Range("B3").Value = "password1"
$ws.Range("C3").Value = 125

# Column D, rows 2-3
$ws.Range("D2").Value = "dev"
$ws.Range("D3").Value = "analyst"

# Row 4
$ws.Range("A4").Value = "username2"
$ws.Range("B4").Value = "password2"
$ws.Range("C4").Value = 126
$ws.Range("D4").Value = "BA"

# E1 becomes a mailto hyperlink (applies the built-in "Hyperlink" style).
$ws.Hyperlinks.Add($ws.Range("E1"), "mailto:rahulriteshqa@gmail.com")

# Add the new, empty "writeData" sheet right after "testData".
$newSheet = $wb.Worksheets.Add($null, $ws)
$newSheet.Name = "writeData"
[void]$newSheet.Range("C16").Select()

# Restore the selection/active sheet back onto "testData".
[void]$ws.Range("D8").Select()
